$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: date bumped by 1 day (45310 -> 45311, i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# D30: price update 570 -> 338
$ws.Range("D30").Value = 338

# D31: price update 690 -> 405
$ws.Range("D31").Value = 405
